# Slide 11, "Espaço Reservado para Texto 5" placeholder: the paragraph
# "Você vai encontrar o pseudocódigo disponível no roteiro." is removed,
# and the following paragraph ("Compare a solução obtida com as versões
# anteriores. ") takes its place, split into two runs: "Compare " and
# "a solução obtida com as versões anteriores. ".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(11)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# Drop the "Compare a solução obtida com as versões anteriores. " paragraph
# (5th paragraph) entirely -- its text will live on, merged into the
# paragraph above.
$oldPara = $tr.Paragraphs(5, 1)
$oldPara.Delete()

# The paragraph that used to read "Você vai encontrar o pseudocódigo
# disponível no roteiro." becomes "Compare a solução obtida com as
# versões anteriores. "
$tr = $sh.TextFrame.TextRange
$targetPara = $tr.Paragraphs(4, 1)
$targetPara.Text = "Compare a solução obtida com as versões anteriores. "

# Split the merged text into two runs: "Compare " + the remainder, so the
# run boundary matches the source edit.
$tr = $sh.TextFrame.TextRange
$targetPara = $tr.Paragraphs(4, 1)
$firstRun = $targetPara.Characters(1, 8)
$firstRun.Text = "Compare "
